$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab and update the "through" date label (May 27 -> May 28)
$ws.Name = "Through 2022-05-28"
$ws.Range("B1").Value = "May 2022 (through May 28)"

# Update / add cell values per neighborhood row
$ws.Range("Q4").Value = 1      # Humboldt Park, May 2019
$ws.Range("AF4").Value = 3     # Humboldt Park, May 2016

$ws.Range("AK5").Value = 2     # Garfield Park, May 2015

$ws.Range("B7").Value = 6      # North Lawndale, May 2022
$ws.Range("L7").Value = 3      # North Lawndale, May 2020

$ws.Range("L11").Value = 4     # Roseland, May 2020

$ws.Range("B15").Value = 3     # Lake View, May 2022

$ws.Range("B25").Value = 6     # Auburn Gresham, May 2022
$ws.Range("L25").Value = 2     # Auburn Gresham, May 2020

$ws.Range("V28").Value = 3     # West Town, May 2018

$ws.Range("L29").Value = 2     # West Pullman, May 2020

$ws.Range("G33").Value = 1     # Avalon Park, May 2021

$ws.Range("B38").Value = 4     # Douglas, May 2022

$ws.Range("AF39").Value = 2    # New City, May 2016

$ws.Range("V45").Value = 1     # Logan Square, May 2018

$ws.Range("AA46").Value = 1    # Little Village, May 2017

$ws.Range("AA51").Value = 1    # Ashburn, May 2017

$ws.Range("Q80").Value = 1     # Oakland, May 2019

$ws.Range("L82").Value = 1     # Portage Park, May 2020
